$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old 2000-2009 rows (rows 2 through 10). This shifts the
# 2010/2011/2012 rows (formerly 11,12,13) up to rows 2,3,4. Those rows
# already have an empty D cell, so nothing further is needed there.
$ws.Range("A2:F10").EntireRow.Delete()
